# Insert a new "CU" column before the existing "nom" (CU) / "url_produit" (CV)
# columns, shifting them one column to the right (-> CV / CW), then populate
# the newly freed CU column with a fresh price-history snapshot:
#   - CU1 (header row) gets the new timestamp "2026-02-01 06:35:39"
#   - CU2:CU80 get the same price value currently held in CT2:CT80
#   - CU81:CU206 stay empty, because CT81:CT206 are already empty there

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift "nom"/"url_produit" (and everything to their right) one column over.
$ws.Columns("CU").Insert()

# New header cell for the freshly inserted column.
$ws.Range("CU1").Value = "2026-02-01 06:35:39"

# Copy this run's price snapshot (column CT, the previous run) into the new
# CU column for every row that actually has a price.
for ($r = 2; $r -le 80; $r++) {
    $price = $ws.Cells.Item($r, 98).Value()
    $ws.Cells.Item($r, 99).Value = $price
}

Write-Output "done"
